$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the answer text in B3: swap order of S: and D: clauses
$ws.Range("B3").Value = "S: Nem, nincs lehetőségem nyugodt körülmények között elvégezni a feladatot.    D: Igen, körülményeim nyugodtak, az internet kapcsolat stabil, a feladatokra tudok szánni 60 percet."

# Bold the font for B3 (new cellXf has applyFont=true, larger/bigger look - "fx cross bigger")
$ws.Range("B3").Font.Bold = $true

# Change selection to B3
$ws.Range("B3").Select()
